$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new (blank) column before column N ("Outstanding -> Due / In Advance / Late"
# split gets an extra spacer column). Everything from N onward (old N,O,P) shifts one
# column to the right (-> O,P,Q).
$ws.Columns("N").Insert()

# The newly inserted column keeps a similar custom width to its neighbours
# (closest reproducible value to the source width of 11.140625 characters).
$ws.Columns("N").ColumnWidth = 10.25

# Make "Repayment Schedule" the active sheet (was "Transactions"), and leave the
# selection where the author last clicked.
$ws.Activate()
$ws.Range("K21").Select()
